$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 7838
$ws1.Range("F14").Value = 29
$ws1.Range("F17").Value = 205
$ws1.Range("F18").Value = 804

# Sheet "全部类型" (All Types) - aggregate view of the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 7838
$ws4.Range("F15").Value = 29
$ws4.Range("F18").Value = 205
$ws4.Range("F19").Value = 804
